$d = $word.ActiveDocument

# Step 1: Reword the first list item's opening sentence.
$d.Content.Find.Execute("Read the", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "Setup Google maps API key and SF Area by lat and long", 2) | Out-Null

# Step 2: Drop the old _GoBack bookmark; we'll recreate it in the right spot
# once the paragraph has been split.
$d.Bookmarks.Item("_GoBack").Delete()

# Step 3: Split " .csv file" off into its own list item (new paragraph),
# re-using the same list/paragraph formatting.
$rng = $d.Content
$rng.Find.Execute(" .csv file", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.InsertParagraphBefore()

# Step 4: Normalize the new paragraph's run text to a single clean run.
$d.Content.Find.Execute(" .csv file", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Read the .csv file", 2) | Out-Null

# Step 5: Recreate the _GoBack bookmark, collapsed, right at the end of the
# "Setup Google..." paragraph (before its paragraph mark). A temporary marker
# character is used so the insertion point isn't the very last position in
# the paragraph when the bookmark is created (which would otherwise anchor
# it to the start of the following paragraph instead).
$p1 = $d.Paragraphs(4)
$pos = $p1.Range.End - 1

$d.Range($pos, $pos).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos)) | Out-Null
$d.Range($pos, $pos + 1).Delete()
